$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while forcing text storage so that
# numeric-looking strings (e.g. "1.00", "0.999") are not auto-converted
# to numbers by Excel, and restore the cell style back to default ("Normal")
# afterward so no stray number-format style is left behind.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Rows 42 and 43: coin name/link swap places, with updated price/volume data ---
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.965"
$ws.Range("E42").Value = "  -1.33%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  -0.04%  "

# --- Price (D) and Volume(1h) (E) updates for remaining rows ---
Set-TextValue "D2" "69.060.03"
$ws.Range("E2").Value = "  +0.30%  "
Set-TextValue "D3" "3.777.22"
$ws.Range("E3").Value = "  -1.01%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue "D5" "628.79"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E6").Value = "  +0.19%  "
Set-TextValue "D7" "3.774.07"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.04%  "
Set-TextValue "D9" "0.520"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  +0.65%  "
Set-TextValue "D12" "6.78"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("E13").Value = "  -4.81%  "
Set-TextValue "D14" "34.84"
$ws.Range("E14").Value = "  -2.83%  "
Set-TextValue "D15" "4.409.28"
$ws.Range("E15").Value = "  -1.06%  "
Set-TextValue "D16" "3.769.90"
$ws.Range("E16").Value = "  -1.35%  "
Set-TextValue "D17" "69.018.26"
$ws.Range("E17").Value = "  +0.21%  "
Set-TextValue "D18" "17.66"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("E19").Value = "  +0.10%  "
Set-TextValue "D20" "7.00"
$ws.Range("E20").Value = "  -1.99%  "
Set-TextValue "D21" "469.01"
$ws.Range("E21").Value = "  +0.81%  "
Set-TextValue "D22" "9.50"
$ws.Range("E22").Value = "  -1.88%  "
Set-TextValue "D23" "0.702"
$ws.Range("E23").Value = "  -0.84%  "
Set-TextValue "D24" "82.08"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  -8.38%  "
Set-TextValue "D26" "12.14"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  -1.48%  "
Set-TextValue "D28" "10.09"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +0.04%  "
Set-TextValue "D30" "3.923.74"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  -0.14%  "
Set-TextValue "D33" "7.08"
$ws.Range("E33").Value = "  -3.28%  "
Set-TextValue "D34" "0.177"
$ws.Range("E34").Value = "  +19.47%  "
Set-TextValue "D35" "28.42"
$ws.Range("E35").Value = "  -2.50%  "
Set-TextValue "D36" "0.999"
$ws.Range("E36").Value = "  -0.28%  "
Set-TextValue "D37" "3.726.33"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("E39").Value = "  -1.04%  "
Set-TextValue "D40" "3.26"
$ws.Range("E40").Value = "  -3.82%  "
Set-TextValue "D41" "5.78"
$ws.Range("E41").Value = "  -2.34%  "
Set-TextValue "D45" "1.99"
$ws.Range("E45").Value = "  +5.59%  "
Set-TextValue "D46" "155.97"
$ws.Range("E46").Value = "  +0.87%  "
Set-TextValue "D47" "43.49"
$ws.Range("E47").Value = "  +1.58%  "
Set-TextValue "D48" "46.98"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("E50").Value = "  -2.62%  "
Set-TextValue "D51" "8.34"
$ws.Range("E51").Value = "  -1.27%  "
